$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows before the current row 351, shifting the
# existing rows 351-381 down to 353-383 (dimension grows to A1:R383).
$ws.Rows("351:352").Insert()

# --- New row 351 ---
$ws.Cells.Item(351, 1).Value2 = 7
$ws.Cells.Item(351, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(351, 3).Value2 = "Ñuble"
$ws.Cells.Item(351, 4).Value2 = 45106
$ws.Cells.Item(351, 5).Value2 = 16
$ws.Cells.Item(351, 6).Value2 = 100112017
$ws.Cells.Item(351, 7).Value2 = "Apio"
$ws.Cells.Item(351, 8).Value2 = "Americana (o)"
$ws.Cells.Item(351, 9).Value2 = "Primera"
$ws.Cells.Item(351, 10).Value2 = 150
$ws.Cells.Item(351, 11).Value2 = 7000
$ws.Cells.Item(351, 12).Value2 = 7000
$ws.Cells.Item(351, 13).Value2 = 7000
$ws.Cells.Item(351, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(351, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(351, 16).Value2 = 1167
$ws.Cells.Item(351, 17).Value2 = 6
$ws.Cells.Item(351, 18).Value2 = "Hortaliza"

# --- New row 352 ---
$ws.Cells.Item(352, 1).Value2 = 7
$ws.Cells.Item(352, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(352, 3).Value2 = "Ñuble"
$ws.Cells.Item(352, 4).Value2 = 45106
$ws.Cells.Item(352, 5).Value2 = 16
$ws.Cells.Item(352, 6).Value2 = 100112017
$ws.Cells.Item(352, 7).Value2 = "Apio"
$ws.Cells.Item(352, 8).Value2 = "Americana (o)"
$ws.Cells.Item(352, 9).Value2 = "Segunda"
$ws.Cells.Item(352, 10).Value2 = 180
$ws.Cells.Item(352, 11).Value2 = 5000
$ws.Cells.Item(352, 12).Value2 = 6000
$ws.Cells.Item(352, 13).Value2 = 5556
$ws.Cells.Item(352, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(352, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(352, 16).Value2 = 926
$ws.Cells.Item(352, 17).Value2 = 6
$ws.Cells.Item(352, 18).Value2 = "Hortaliza"
